# Applies the "Player Info" sheet addition + MATCH_CARD_LINK -> MATCH_CODE
# rewrite (URL -> bare match-code value) to ODI Batting / ODI Bowling.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a new "Player Info" worksheet as the very first sheet in the book.
#    We clone "ODI Batting" (Worksheet.Copy) rather than Worksheets.Add so
#    the header row naturally inherits the exact same header style used on
#    the other sheets (bold / bordered / centered-top) instead of having to
#    rebuild it cell-by-cell, which only manages to approximate it.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Copy($battingSheet)
$playerInfo = $wb.Worksheets.Item(1)
$playerInfo.Name = "Player Info"

# Drop the cloned ODI Batting content; ClearContents keeps per-cell styles
# (unlike Clear, which wipes formatting too) so the styled A1:D1 header
# cells are left intact for us to overwrite with the new header text.
$playerInfo.Cells.ClearContents()

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# ID is stored as text (e.g. "4273"), same as every other column on this
# sheet, so force text storage instead of letting it coerce to a number.
$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "4273"
$idCell.ClearFormats()

$playerInfo.Range("B2").Value = "Beuran Eric Hendricks"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Fast Medium"

# Remove the leftover cells copied in from ODI Batting that aren't part of
# the new 4-column x 2-row Player Info sheet (the extra header-styled
# columns E:J on row 1, and old data rows 2-9).
$playerInfo.Range("E1:J1").Clear()
$playerInfo.Range("A3:J9").Clear()

# ---------------------------------------------------------------------------
# 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE, URL values -> bare codes.
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingUsed = $batting.UsedRange
$battingRows = $battingUsed.Rows.Count
for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        $code = $val -replace '^.*MatchCode=', ''
        # Force the replacement to be stored as text (not auto-coerced to a
        # number) the same way the source value was a text/inlineStr cell.
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------------
# 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE, URL values -> bare codes.
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingUsed = $bowling.UsedRange
$bowlingRows = $bowlingUsed.Rows.Count
for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        $code = $val -replace '^.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.ClearFormats()
    }
}
